$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M")

# The "period" column (A2:A4) used to hold text like "2020M01"; switch it
# to real date values so downstream consumers can treat it as a date.
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Value = 43831   # 2020-01-01

# Reuse the same formatting (and therefore the same cell style) for the
# rest of the column instead of re-deriving it cell by cell.
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$ws.Range("A3").Value = 43862   # 2020-02-01
$ws.Range("A4").Value = 43892   # 2020-03-02

# Reflect the edited range in the sheet's selection, then restore the
# workbook's originally-active sheet ("Q").
$null = $ws.Select()
$null = $ws.Range("A2:A4").Select()

$null = $wb.Worksheets.Item("Q").Select()
